$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added clock in times: each shift-role label ("1: Server", "H/G", "Bar", "Expo", ...)
# now carries the clock-in time as a prefix, e.g. "1: Server" -> "10:00: 1",
# "H/G" -> "11:00: H/G" / "4:00: H/G", "Bar"/"Expo" -> "10:00: ..." / "4:00: ...".
# Write the exact post-edit label into every cell whose text changed.

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("D2").Value = "10:00: Bar"
$ws.Range("G2").Value = "4:00: Bar"
$ws.Range("I2").Value = "-"
$ws.Range("L2").Value = "-"
$ws.Range("N2").Value = "-"
# Row 3
$ws.Range("B3").Value = "10:00: Bar"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "10:00: Bar"
$ws.Range("K3").Value = "-"
$ws.Range("M3").Value = "4:00: Bar"
$ws.Range("N3").Value = "10:00: Bar"
# Row 4
$ws.Range("C4").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "4:00: Bar"
$ws.Range("I4").Value = "4:00: Bar"
$ws.Range("J4").Value = "-"
$ws.Range("K4").Value = "4:00: Bar"
$ws.Range("M4").Value = "4:00: Bar"
# Row 5
$ws.Range("C5").Value = "4:00: Bar"
$ws.Range("E5").Value = "4:00: Bar"
$ws.Range("J5").Value = "10:00: Bar"
$ws.Range("K5").Value = "4:00: Bar"
$ws.Range("N5").Value = "10:00: Bar"
$ws.Range("O5").Value = "4:00: Bar"
# Row 6
$ws.Range("H6").Value = "10:00: Bar"
$ws.Range("I6").Value = "4:00: Bar"
$ws.Range("N6").Value = "-"
# Row 7
$ws.Range("B7").Value = "10:00: Bar"
$ws.Range("C7").Value = "4:00: Bar"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "4:00: Bar"
$ws.Range("G7").Value = "-"
$ws.Range("I7").Value = "-"
$ws.Range("L7").Value = "10:00: Bar"
$ws.Range("M7").Value = "-"
$ws.Range("O7").Value = "4:00: Bar"
# Row 8
$ws.Range("B8").Value = "10:00: Expo"
$ws.Range("D8").Value = "10:00: 1"
$ws.Range("E8").Value = "5:00: 1"
$ws.Range("F8").Value = "10:00: Expo"
$ws.Range("G8").Value = "5:00: 1"
$ws.Range("H8").Value = "10:00: 1"
$ws.Range("I8").Value = "4:00: 3"
$ws.Range("J8").Value = "-"
$ws.Range("K8").Value = "4:00: 3"
$ws.Range("L8").Value = "10:00: 1"
$ws.Range("M8").Value = "4:00: Expo"
$ws.Range("N8").Value = "10:00: 1"
$ws.Range("O8").Value = "4:00: Expo"
# Row 9
$ws.Range("D9").Value = "-"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "4:00: 3"
$ws.Range("I9").Value = "5:00: 2"
$ws.Range("J9").Value = "10:00: 3"
$ws.Range("L9").Value = "-"
$ws.Range("N9").Value = "-"
# Row 10
$ws.Range("C10").Value = "-"
$ws.Range("E10").Value = "-"
$ws.Range("K10").Value = "-"
$ws.Range("M10").Value = "4:00: 3"
# Row 11
$ws.Range("B11").Value = "10:00: 3"
$ws.Range("H11").Value = "-"
$ws.Range("I11").Value = "4:00: 5"
$ws.Range("L11").Value = "10:00: 5"
# Row 12
$ws.Range("B12").Value = "-"
$ws.Range("G12").Value = "-"
$ws.Range("M12").Value = "4:00: 4"
$ws.Range("N12").Value = "10:00: 5"
# Row 13
$ws.Range("B13").Value = "10:00: 4"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "4:00: 4"
$ws.Range("L13").Value = "-"
# Row 14
$ws.Range("C14").Value = "-"
$ws.Range("E14").Value = "5:00: 2"
$ws.Range("F14").Value = "10:00: 1"
$ws.Range("G14").Value = "-"
$ws.Range("K14").Value = "4:00: 5"
$ws.Range("O14").Value = "4:00: 5"
# Row 15
$ws.Range("B15").Value = "10:00: 1"
$ws.Range("F15").Value = "10:00: 2"
$ws.Range("G15").Value = "5:00: 2"
$ws.Range("I15").Value = "5:00: 1"
$ws.Range("J15").Value = "10:00: 1"
# Row 16
$ws.Range("C16").Value = "4:00: 2"
$ws.Range("D16").Value = "10:00: 2"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "10:00: 3"
$ws.Range("H16").Value = "-"
$ws.Range("I16").Value = "-"
$ws.Range("J16").Value = "-"
$ws.Range("K16").Value = "5:00: 2"
$ws.Range("L16").Value = "-"
$ws.Range("M16").Value = "4:00: 1"
$ws.Range("O16").Value = "4:00: 1"
# Row 17
$ws.Range("B17").Value = "-"
$ws.Range("C17").Value = "4:00: 3"
$ws.Range("F17").Value = "10:00: 4"
$ws.Range("H17").Value = "10:00: 3"
$ws.Range("I17").Value = "-"
$ws.Range("J17").Value = "10:00: 4"
$ws.Range("L17").Value = "10:00: 2"
$ws.Range("N17").Value = "10:00: 2"
# Row 18
$ws.Range("B18").Value = "10:00: 5"
$ws.Range("L18").Value = "10:00: 4"
# Row 19
$ws.Range("C19").Value = "4:00: 4"
$ws.Range("K19").Value = "4:00: 4"
# Row 20
$ws.Range("B20").Value = "10:00: 2"
$ws.Range("G20").Value = "-"
$ws.Range("H20").Value = "10:00: 4"
# Row 21
$ws.Range("C21").Value = "-"
$ws.Range("G21").Value = "4:00: 5"
$ws.Range("L21").Value = "-"
$ws.Range("N21").Value = "-"
$ws.Range("O21").Value = "4:00: 4"
# Row 22
$ws.Range("C22").Value = "4:00: 1"
$ws.Range("G22").Value = "-"
$ws.Range("H22").Value = "-"
$ws.Range("I22").Value = "-"
$ws.Range("J22").Value = "10:00: 2"
$ws.Range("K22").Value = "5:00: 1"
$ws.Range("L22").Value = "10:00: 3"
$ws.Range("N22").Value = "-"
$ws.Range("O22").Value = "4:00: 3"
# Row 23
$ws.Range("C23").Value = "4:00: 5"
$ws.Range("D23").Value = "10:00: 3"
$ws.Range("E23").Value = "4:00: 5"
$ws.Range("H23").Value = "10:00: 2"
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("M23").Value = "4:00: 2"
$ws.Range("N23").Value = "-"
$ws.Range("O23").Value = "-"
# Row 24
$ws.Range("D24").Value = "10:00: 4"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("H24").Value = "-"
$ws.Range("K24").Value = "-"
$ws.Range("M24").Value = "-"
$ws.Range("N24").Value = "10:00: 4"
$ws.Range("O24").Value = "4:00: 2"
# Row 25
$ws.Range("B25").Value = "-"
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "4:00: 3"
$ws.Range("G25").Value = "4:00: 4"
$ws.Range("I25").Value = "4:00: 4"
$ws.Range("M25").Value = "4:00: 5"
$ws.Range("N25").Value = "10:00: 3"
# Row 26
$ws.Range("B26").Value = "-"
$ws.Range("D26").Value = "-"
$ws.Range("H26").Value = "-"
$ws.Range("L26").Value = "11:00: H/G"
# Row 27
$ws.Range("G27").Value = "-"
$ws.Range("H27").Value = "11:00: H/G"
$ws.Range("I27").Value = "-"
$ws.Range("J27").Value = "11:00: H/G"
# Row 28
$ws.Range("B28").Value = "11:00: H/G"
$ws.Range("C28").Value = "-"
$ws.Range("N28").Value = "11:00: H/G"
# Row 29
$ws.Range("C29").Value = "4:00: H/G"
$ws.Range("D29").Value = "11:00: H/G"
$ws.Range("E29").Value = "-"
$ws.Range("I29").Value = "4:00: H/G"
$ws.Range("J29").Value = "-"
$ws.Range("O29").Value = "4:00: H/G"
# Row 30
$ws.Range("E30").Value = "4:00: H/G"
$ws.Range("F30").Value = "11:00: H/G"
$ws.Range("G30").Value = "4:00: H/G"
$ws.Range("K30").Value = "4:00: H/G"
$ws.Range("M30").Value = "4:00: H/G"
$ws.Range("O30").Value = "-"
# Row 31
$ws.Range("F31").Value = "-"
$ws.Range("K31").Value = "-"
$ws.Range("L31").Value = "-"
# Row 32
$ws.Range("B32").Value = "-"
$ws.Range("C32").Value = "4:00: Expo"
$ws.Range("D32").Value = "10:00: Expo"
$ws.Range("E32").Value = "4:00: Expo"
$ws.Range("F32").Value = "-"
$ws.Range("G32").Value = "4:00: Expo"
$ws.Range("H32").Value = "10:00: Expo"
$ws.Range("I32").Value = "4:00: Expo"
$ws.Range("J32").Value = "10:00: Expo"
$ws.Range("K32").Value = "4:00: Expo"
$ws.Range("L32").Value = "10:00: Expo"
$ws.Range("M32").Value = "4:00: Expo"
$ws.Range("N32").Value = "10:00: Expo"
$ws.Range("O32").Value = "4:00: Expo"
